# Interview prep.docx edit
#
# Content changes (per commit "O(N) reverse polish notation" and the
# supplied OOXML diff):
#
#   1. "evaluate-reverse-polish-notation" paragraph: the trailing run
#      (a single space after the hyperlink) gets "Done" appended as its
#      own new run -> " " + "Done".
#
#   2. "min-stack" paragraph: the note "(done, just upload it)" is
#      replaced with "done" (kept as its own run, after the existing
#      leading-space run).
#
# (All the other hunks in the source diff are pure Word-grammar-checker
# proofErr bookkeeping / incidental run splitting around unchanged text
# - "duplicate", "Done", "5"+"3", etc. - with no visible-text change, so
# they are not re-created here.)

$d = $word.ActiveDocument

function Get-ParagraphByContains($needle) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# Locate the absolute Range of $needle (a substring of the *visible*
# text) inside paragraph $para, walking character-by-character so that
# hidden/field characters in the underlying Range don't throw off the
# offsets.
function Find-SubrangeByText($para, $needle) {
    $r = $para.Range
    $visible = ""
    $posMap = @()
    for ($i = $r.Start; $i -lt $r.End; $i++) {
        $ch = $d.Range($i, $i + 1).Text
        if ($ch -ne "") {
            $visible += $ch
            $posMap += $i
        }
    }
    $idx = $visible.IndexOf($needle)
    if ($idx -lt 0) { return $null }
    $startPos = $posMap[$idx]
    $endPos = $posMap[$idx + $needle.Length - 1] + 1
    return $d.Range($startPos, $endPos)
}

# --- 1. Evaluate Reverse Polish Notation -> mark Done -----------------
$p1 = Get-ParagraphByContains("evaluate-reverse-polish-notation")
if ($p1 -ne $null) {
    $end1 = $p1.Range.End
    # Collapse to just before the paragraph mark, then append "Done" as
    # a brand-new run (so it doesn't merge into the preceding run).
    $insertAt1 = $d.Range($end1 - 1, $end1 - 1)
    $insertAt1.InsertAfter("Done")
}

# --- 2. Min Stack -> "(done, just upload it)" becomes "done" ----------
$p2 = Get-ParagraphByContains("min-stack/")
if ($p2 -ne $null) {
    $sub2 = Find-SubrangeByText $p2 "(done, just upload it)"
    if ($sub2 -ne $null) {
        $insertPos2 = $sub2.Start
        $sub2.Text = ""
        $insertAt2 = $d.Range($insertPos2, $insertPos2)
        $insertAt2.InsertAfter("done")
    }
}

Write-Host "Done applying edits."
